# "testing 2050 scenario with batteries"
# Applies the parameter edits made on the "Coupling Parameters" sheet:
#  - Start Year:              2020 -> 2050
#  - End Year:                2050 -> 2055
#  - Power_plants_from_year:  2020 -> 2050
#  - fix_price_year:          2020 -> 2050
#  - fix_demand_to_initial_year: FALSE -> TRUE
#  - targetinvestment_per_year:  TRUE -> FALSE
#  - extend the "yearly_CO2_prices" note text
#  - taller row for the start_tick_dismantling note (21.5 -> 31.5)
#  - move the active selection from B26 to C22

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Start Year: 2020 -> 2050
$ws.Range("B2").Value = 2050

# End Year: 2050 -> 2055
$ws.Range("B3").Value = 2055

# Power_plants_from_year: 2020 -> 2050
$ws.Range("B4").Value = 2050

# Make the start_tick_dismantling description row taller
$ws.Rows.Item(14).RowHeight = 31.5

# Update the yearly_CO2_prices note to explain the fixed price behaviour
$ws.Range("C16").Value = "so far this is only for NL. If False then the price is fixed to the fix_price_year"

# fix_price_year: 2020 -> 2050
$ws.Range("B18").Value = 2050

# fix_demand_to_initial_year: FALSE -> TRUE
$ws.Range("B19").Value = $true

# targetinvestment_per_year: TRUE -> FALSE
$ws.Range("B25").Value = $false

# Move the selection to C22 (as last left by the author)
$ws.Range("C22").Select()
